# Update plots for each sample.
#
# The marker CYP2D6_002 / CYP2D6_10B (sample S1) and marker CYP2D6_011 /
# CYP2D6_4 (sample S2) were re-analyzed: the mutant peak-call for S1 and
# the wildtype peak-call for S2 that previously could not be detected are
# now found, which updates the peak_table, allele_table and marker_table
# sheets accordingly.

$wb = $excel.ActiveWorkbook

# --- peak_table: adjust the binning window heights/max used for detection ---
$ws1 = $wb.Worksheets.Item("peak_table")
$ws1.Range("N3").Value = 500     # S1 / CYP2D6_002 (CYP2D6_10B): w_height 1000 -> 500
$ws1.Range("G12").Value = 33     # S2 / CYP2D6_011 (CYP2D6_4): w_max 32 -> 33

# --- allele_table: the two previously-undetected peaks are now detected ---
$ws2 = $wb.Worksheets.Item("allele_table")

# Row 4: S1 / CYP2D6_002 (CYP2D6_10B), Forward, mutant base "C"
$ws2.Range("K4").Value = 500
$ws2.Range("M4").Value = $true
$ws2.Range("N4").Value = 48
$ws2.Range("O4").Value = 33.25
$ws2.Range("P4").Value = 603
$ws2.Range("Q4").Value = "ok"
$ws2.Range("R4").Value = ""

# Row 22: S2 / CYP2D6_011 (CYP2D6_4), Forward, wildtype base "G"
$ws2.Range("J22").Value = 33
$ws2.Range("M22").Value = $true
$ws2.Range("N22").Value = 53
$ws2.Range("O22").Value = 31.84
$ws2.Range("P22").Value = 5718
$ws2.Range("Q22").Value = "ok"
$ws2.Range("R22").Value = ""

# --- marker_table: resulting genotype/phenotype calls ---
$ws3 = $wb.Worksheets.Item("marker_table")

# Row 3: S1 / CYP2D6_002 (CYP2D6_10B) now calls heterozygous C/T
$ws3.Range("G3").Value = "CT"
$ws3.Range("H3").Value = "heterozygous"

# Row 12: S2 / CYP2D6_011 (CYP2D6_4) now calls wildtype G/G
$ws3.Range("G12").Value = "GG"
$ws3.Range("H12").Value = "wildtype"
